# Convert some images to code
# - Reposition the 4 remaining diagram shapes (they used to sit "behind"/
#   next to the image that showed the code; now that the code textbox is
#   gone, the diagram shapes shift down-and-right to the spot the image
#   used to occupy).
# - Remove the "Content Placeholder 3" shape, which held a screenshot-style
#   rendering of the Account/Guarantor code (replaced elsewhere by real
#   text, per the commit message "Convert some images to code").

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Values below are expressed in points (1 pt = 12700 EMU) but chosen so
# that the internal float32 round-trip used by this host lands on the
# exact target EMU offsets from the canonical OOXML.

$rect3 = $s.Shapes.Item("Rectangle 3")
$rect3.Left = 227.36009216308594
$rect3.Top = 205.60401916503906

$rect4 = $s.Shapes.Item("Rectangle 4")
$rect4.Left = 401.33984375
$rect4.Top = 208.60401916503906

$elbow5 = $s.Shapes.Item("Elbow Connector 5")
$elbow5.Left = 327.07049560546875
$elbow5.Top = 232.54150390625

$rect6 = $s.Shapes.Item("Rectangle 6")
$rect6.Left = 371.33984375
$rect6.Top = 196.60401916503906

# Remove the code-screenshot placeholder entirely.
$codeBox = $s.Shapes.Item("Content Placeholder 3")
$codeBox.Delete()
